$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1023
$ws.Range("I32").Value = 900
$ws.Range("J32").Value = 1064
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 1064
$ws.Range("M32").Value = -574
$ws.Range("N32").Value = -1716

$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1325
$ws.Range("N40").ClearContents()

$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11944

$ws.Range("H51").Value = 5889.5
$ws.Range("J51").Value = 6702.385
$ws.Range("L51").Value = 6702.385
$ws.Range("N51").Value = -7670.385

$ws.Range("H70").Value = 2835.889
$ws.Range("I70").Value = 3250
$ws.Range("J70").Value = 2504.6
$ws.Range("K70").Value = 9750
$ws.Range("L70").Value = 7513.799999999999
$ws.Range("M70").Value = -9480
$ws.Range("N70").Value = -8053.799999999999

$ws.Range("H73").Value = 2835.889
$ws.Range("I73").Value = 3250
$ws.Range("J73").Value = 2504.6
$ws.Range("K73").Value = 9750
$ws.Range("L73").Value = 7513.799999999999
$ws.Range("M73").Value = -8814
$ws.Range("N73").Value = -9385.799999999999

$ws.Range("H135").Value = 18321.346
$ws.Range("I135").Value = 22809.934
$ws.Range("K135").Value = 205289.406
$ws.Range("M135").Value = -202754.406
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 3336539
$ws.Range("I137").Value = 4350216
$ws.Range("J137").Value = 5885.5713
$ws.Range("K137").Value = 13050648
$ws.Range("L137").Value = 17656.7139
$ws.Range("M137").Value = -13048098
$ws.Range("N137").Value = -22756.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2034.3
$ws.Range("I2").Value = 1979.1428
$ws.Range("K2").Value = 1979.1428
$ws.Range("M2").Value = -1866.1428
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 2123.19
$ws.Range("I32").Value = 1723.8846
$ws.Range("J32").Value = 3538.9092
$ws.Range("K32").Value = 1723.8846
$ws.Range("L32").Value = 3538.9092
$ws.Range("M32").Value = -1436.8846
$ws.Range("N32").Value = -4112.9092

$ws.Range("H74").Value = 6463845.5
$ws.Range("I74").Value = 7845367.5
$ws.Range("J74").Value = 148316.28
$ws.Range("K74").Value = 7845367.5
$ws.Range("L74").Value = 148316.28
$ws.Range("M74").Value = -7844493.5
$ws.Range("N74").Value = -150064.28

$ws.Range("H77").Value = 6463845.5
$ws.Range("I77").Value = 7845367.5
$ws.Range("J77").Value = 148316.28
$ws.Range("K77").Value = 39226837.5
$ws.Range("L77").Value = 741581.4
$ws.Range("M77").Value = -39222469.5
$ws.Range("N77").Value = -750317.4

$ws.Range("H116").Value = 2034.3
$ws.Range("I116").Value = 1979.1428
$ws.Range("K116").Value = 1979.1428
$ws.Range("M116").Value = 314.8571999999999
$ws.Range("N116").ClearContents()

$ws.Range("H133").Value = 29130.5
$ws.Range("J133").Value = 29130.5
$ws.Range("L133").Value = 29130.5
$ws.Range("N133").Value = -34190.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2034.3
$ws.Range("I3").Value = 1979.1428
$ws.Range("K3").Value = 1979.1428
$ws.Range("M3").Value = -1865.1428
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 340
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 366.66666
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 366.66666
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -592.66666

$ws.Range("H4").Value = 3354666.8
$ws.Range("J4").Value = 3354666.8
$ws.Range("L4").Value = 3354666.8
$ws.Range("N4").Value = -3354890.8

$ws.Range("H10").Value = 183.57143
$ws.Range("I10").Value = 183.57143
$ws.Range("K10").Value = 183.57143
$ws.Range("M10").Value = -44.57142999999999

$ws.Range("H11").Value = 800
$ws.Range("J11").Value = 800
$ws.Range("L11").Value = 800
$ws.Range("N11").Value = -1080

$ws.Range("H12").Value = 11498.75
$ws.Range("I12").Value = 930
$ws.Range("J12").Value = 17840
$ws.Range("K12").Value = 930
$ws.Range("L12").Value = 17840
$ws.Range("M12").Value = -760
$ws.Range("N12").Value = -18180

$ws.Range("H19").Value = 619.6667
$ws.Range("I19").Value = 339.66666
$ws.Range("J19").Value = 899.6667
$ws.Range("K19").Value = 339.66666
$ws.Range("L19").Value = 899.6667
$ws.Range("M19").Value = -169.66666
$ws.Range("N19").Value = -1239.6667

$ws.Range("H24").Value = 619.6667
$ws.Range("I24").Value = 339.66666
$ws.Range("J24").Value = 899.6667
$ws.Range("K24").Value = 339.66666
$ws.Range("L24").Value = 899.6667
$ws.Range("M24").Value = -169.66666
$ws.Range("N24").Value = -1239.6667

$ws.Range("H86").Value = 2129.0334
$ws.Range("I86").Value = 1963
$ws.Range("J86").Value = 2585.625
$ws.Range("K86").Value = 1963
$ws.Range("L86").Value = 2585.625
$ws.Range("M86").Value = -840
$ws.Range("N86").Value = -4831.625

$ws.Range("H89").Value = 2129.0334
$ws.Range("I89").Value = 1963
$ws.Range("J89").Value = 2585.625
$ws.Range("K89").Value = 9815
$ws.Range("L89").Value = 12928.125
$ws.Range("M89").Value = -4199
$ws.Range("N89").Value = -24160.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13514430
$ws.Range("J131").Value = 1047.7167
$ws.Range("L131").Value = 3143.1501
$ws.Range("N131").Value = -13223.1501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5833667.5
$ws.Range("I11").Value = 6692461
$ws.Range("J11").Value = 251510
$ws.Range("K11").Value = 6692461
$ws.Range("L11").Value = 251510
$ws.Range("M11").Value = -6692322
$ws.Range("N11").Value = -251788

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1032.091
$ws.Range("I9").Value = 257.8
$ws.Range("J9").Value = 1677.3334
$ws.Range("K9").Value = 257.8
$ws.Range("L9").Value = 1677.3334
$ws.Range("M9").Value = -33.80000000000001
$ws.Range("N9").Value = -2125.3334

$ws.Range("H10").Value = 7500351
$ws.Range("I10").Value = 15000202
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 15000202
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = -15000062
$ws.Range("N10").Value = -780

$ws.Range("H12").Value = 2000
$ws.Range("I12").Value = 3500
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 3500
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -3330
$ws.Range("N12").Value = -840

$ws.Range("H132").Value = 33660.688
$ws.Range("I132").Value = 15313.553
$ws.Range("J132").Value = 103379.8
$ws.Range("K132").Value = 45940.659
$ws.Range("L132").Value = 310139.4
$ws.Range("M132").Value = -43410.659
$ws.Range("N132").Value = -315199.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 125004360
$ws.Range("I62").Value = 250003710
$ws.Range("J62").Value = 5003
$ws.Range("K62").Value = 250003710
$ws.Range("L62").Value = 5003
$ws.Range("M62").Value = -250003086
$ws.Range("N62").Value = -6251

$ws.Range("H65").Value = 125004360
$ws.Range("I65").Value = 250003710
$ws.Range("J65").Value = 5003
$ws.Range("K65").Value = 1250018550
$ws.Range("L65").Value = 25015
$ws.Range("M65").Value = -1250015430
$ws.Range("N65").Value = -31255

$ws.Range("H136").Value = 41260.22
$ws.Range("I136").Value = 33099.582
$ws.Range("J136").Value = 54574.95
$ws.Range("K136").Value = 99298.74600000001
$ws.Range("L136").Value = 163724.85
$ws.Range("M136").Value = -96748.74600000001
$ws.Range("N136").Value = -168824.85
